$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.518.58'
$ws.Range('E2').Value = '  +3.40%  '
$ws.Range('D3').Value = '2.639.36'
$ws.Range('E3').Value = '  +0.49%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '570.67'
$ws.Range('E5').Value = '  +6.63%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.92'
$ws.Range('E6').Value = '  +2.71%  '
$ws.Range('E7').Value = '  -0.31%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.610'
$ws.Range('E8').Value = '  +7.36%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '6.82'
$ws.Range('E9').Value = '  -2.15%  '
$ws.Range('E10').Value = '  +3.94%  '
$ws.Range('E11').Value = '  +6.42%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.343'
$ws.Range('E12').Value = '  +2.71%  '
$ws.Range('D13').Value = '3.109.48'
$ws.Range('E13').Value = '  +0.54%  '
$ws.Range('D14').Value = '60.522.82'
$ws.Range('E14').Value = '  +3.52%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '21.74'
$ws.Range('E15').Value = '  +4.70%  '
$ws.Range('E16').Value = '  +4.00%  '
$ws.Range('D17').Value = '2.652.97'
$ws.Range('E17').Value = '  +1.60%  '
$ws.Range('E18').Value = '  +3.58%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '345.26'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.43'
$ws.Range('E20').Value = '  +2.84%  '
$ws.Range('E21').Value = '  +2.67%  '
$ws.Range('E22').Value = '  +1.18%  '
$ws.Range('E23').Value = '  +0.01%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '66.88'
$ws.Range('E24').Value = '  +0.92%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.443'
$ws.Range('E25').Value = '  +6.67%  '
$ws.Range('E26').Value = '  +2.12%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.993'
$ws.Range('E27').Value = '  -0.54%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.36'
$ws.Range('E28').Value = '  +3.51%  '
$ws.Range('E29').Value = '  +5.60%  '
$ws.Range('E30').Value = '  -0.12%  '
$ws.Range('E31').Value = '  +4.79%  '
$ws.Range('E32').Value = '  +4.16%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '156.16'
$ws.Range('E33').Value = '  +3.89%  '
$ws.Range('E34').Value = '  +2.31%  '
$ws.Range('E35').Value = '  +5.19%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.915'
$ws.Range('E36').Value = '  +7.70%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.18'
$ws.Range('E37').Value = '  +6.70%  '
$ws.Range('B38').Value = 'Fetch.AI'
$ws.Range('C38').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.910'
$ws.Range('E38').Value = '  +12.31%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '37.65'
$ws.Range('E39').Value = '  +1.22%  '
$ws.Range('E40').Value = '  +7.43%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '307.79'
$ws.Range('E42').Value = '  +2.93%  '
$ws.Range('E43').Value = '  -0.47%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.608'
$ws.Range('E44').Value = '  +2.01%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0981'
$ws.Range('E45').Value = '  +4.75%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0550'
$ws.Range('E46').Value = '  +3.78%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '19.46'
$ws.Range('E47').Value = '  +2.29%  '
$ws.Range('E48').Value = '  -0.14%  '
$ws.Range('E49').Value = '  +5.30%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '125.41'
$ws.Range('E50').Value = '  +11.28%  '
$ws.Range('D51').Value = '1.970.83'
$ws.Range('E51').Value = '  +1.21%  '
